# "Updated LoginPage with ErrMsg"
# Adds a new "ValidLogin" worksheet (after the existing "TestCase1" sheet)
# that holds a small login fixture table: UserName/Password headers with
# an admin/manager row beneath them. The new sheet becomes the active
# (selected) sheet/tab, matching the authored workbook.

$wb = $excel.ActiveWorkbook

# Existing first sheet -> new sheet gets inserted right after it so it
# keeps sheetId="1"/rId1 and the new sheet becomes sheetId="2"/rId2.
$ws1 = $wb.Worksheets.Item("TestCase1")

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "ValidLogin"

# Header row
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"

# Data row
$ws2.Range("A2").Value = "admin"
$ws2.Range("B2").Value = "manager"

# Match the authored selection/active cell and zoom on the new sheet.
$win = $wb.Windows.Item(1)
$win.Zoom = 160
$ws2.Range("B3").Select()
